$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Initial Investment ---
$ws.Range("B2").Value = -1000000

# --- Row 3: Depreciation (C3:L3) ---
$ws.Range("C3:L3").Value = 30000

# --- Row 4: Incoming Payments (C4:L4) ---
$ws.Range("C4:L4").Value = 1400000

# --- Row 5: Outgoing Payments ---
$ws.Range("B5").Value = -350000
$ws.Range("C5:L5").Value = -700000

# --- Row 8: Yearly Net ---
$ws.Range("B8").Value = -1550000
$ws.Range("C8:K8").Value = 730000
$ws.Range("L8").Value = 1070000

# --- Row 9: Present Value ---
$ws.Range("B9").Value = -1550000
$ws.Range("C9").Value = 673431.7343173431
$ws.Range("D9").Value = 621246.9873776227
$ws.Range("E9").Value = 573106.0769166261
$ws.Range("F9").Value = 528695.6429120167
$ws.Range("G9").Value = 487726.6078524138
$ws.Range("H9").Value = 449932.2950668024
$ws.Range("I9").Value = 415066.6928660539
$ws.Range("J9").Value = 382902.8531974667
$ws.Range("K9").Value = 353231.4143888069
$ws.Range("L9").Value = 477629.2945913452

# --- Row 10: Accumulated Present Value ---
$ws.Range("B10").Value = -1550000
$ws.Range("C10").Value = -876568.2656826569
$ws.Range("D10").Value = -255321.2783050342
$ws.Range("E10").Value = 317784.7986115919
$ws.Range("F10").Value = 846480.4415236086
$ws.Range("G10").Value = 1334207.049376023
$ws.Range("H10").Value = 1784139.344442825
$ws.Range("I10").Value = 2199206.037308879
$ws.Range("J10").Value = 2582108.890506345
$ws.Range("K10").Value = 2935340.304895152
$ws.Range("L10").Value = 3412969.599486498

# The values in E10:H10 flipped from negative to positive, so the manual
# red/green highlighting needs to follow (same visual convention already
# used by I10:L10, which carry the "positive" green style). Copy that
# cell format across instead of fabricating a brand-new style.
$ws.Range("I10").Copy()
$ws.Range("E10:H10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 11: Net Present Value ---
$ws.Range("B11").Value = 3412969.599486498
